$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 48; existing rows 48-52 shift down to 49-53.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row with the new team.
$ws.Range("A48").Value = "TEAM LOPES 99"
$ws.Range("B48").Value = 479510
$ws.Range("C48").Value = "https://cartola.globo.com/#!/time/479510"

# Row-insert does not renumber the worksheet's stored hyperlink refs, so
# rebuild the whole Hyperlinks collection for C2:C53 against the correct rows.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "https://cartola.globo.com/", "!/time/117598") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://cartola.globo.com/", "!/time/49355335") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://cartola.globo.com/", "!/time/18346776") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://cartola.globo.com/", "!/time/13913874") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://cartola.globo.com/", "!/time/50988641") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://cartola.globo.com/", "!/time/25748736") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://cartola.globo.com/", "!/time/3851966") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://cartola.globo.com/", "!/time/49243759") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://cartola.globo.com/", "!/time/7017989") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "https://cartola.globo.com/", "!/time/387186") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "https://cartola.globo.com/", "!/time/20696550") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "https://cartola.globo.com/", "!/time/186283") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "https://cartola.globo.com/", "!/time/1863710") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), "https://cartola.globo.com/", "!/time/48279389") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), "https://cartola.globo.com/", "!/time/25311459") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "https://cartola.globo.com/", "!/time/18642587") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), "https://cartola.globo.com/", "!/time/25565675") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), "https://cartola.globo.com/", "!/time/18344271") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C20"), "https://cartola.globo.com/", "!/time/18421230") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C21"), "https://cartola.globo.com/", "!/time/528730") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C22"), "https://cartola.globo.com/", "!/time/24468241") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C23"), "https://cartola.globo.com/", "!/time/13951133") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C24"), "https://cartola.globo.com/", "!/time/1747619") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C25"), "https://cartola.globo.com/", "!/time/51010813") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C26"), "https://cartola.globo.com/", "!/time/44810918") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C27"), "https://cartola.globo.com/", "!/time/20340994") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C28"), "https://cartola.globo.com/", "!/time/4911779") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C29"), "https://cartola.globo.com/", "!/time/19033717") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C30"), "https://cartola.globo.com/", "!/time/25401606") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C31"), "https://cartola.globo.com/", "!/time/30267301") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C32"), "https://cartola.globo.com/", "!/time/3708025") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C33"), "https://cartola.globo.com/", "!/time/14124559") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C34"), "https://cartola.globo.com/", "!/time/48498051") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C35"), "https://cartola.globo.com/", "!/time/25313333") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C36"), "https://cartola.globo.com/", "!/time/9823692") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C37"), "https://cartola.globo.com/", "!/time/3447341") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C38"), "https://cartola.globo.com/", "!/time/18223508") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C39"), "https://cartola.globo.com/", "!/time/5823700") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C40"), "https://cartola.globo.com/", "!/time/29228373") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C41"), "https://cartola.globo.com/", "!/time/25811332") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C42"), "https://cartola.globo.com/", "!/time/1148959") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C43"), "https://cartola.globo.com/", "!/time/13707047") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C44"), "https://cartola.globo.com/", "!/time/4229593") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C45"), "https://cartola.globo.com/", "!/time/28741323") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C46"), "https://cartola.globo.com/", "!/time/49180400") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C47"), "https://cartola.globo.com/", "!/time/212042") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C48"), "https://cartola.globo.com/", "!/time/479510") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C49"), "https://cartola.globo.com/", "!/time/335716") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C50"), "https://cartola.globo.com/", "!/time/1273719") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C51"), "https://cartola.globo.com/", "!/time/3424598") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C52"), "https://cartola.globo.com/", "!/time/2981301") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C53"), "https://cartola.globo.com/", "!/time/14696986") | Out-Null

# Re-adding hyperlinks nudges the cell style; restore the original "Hyperlink" style.
$ws.Range("C2:C53").Style = "Hyperlink"
